$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Seed shared-string label column (B2:B19) top-to-bottom ---
# (matches natural first-use ordering for the label strings)
$ws.Range("B2").Value = 'HKL'
$ws.Range("B3").Value = 'ND Single'
$ws.Range("B4").Value = 'RD Single'
$ws.Range("B5").Value = 'TD Single'
$ws.Range("B6").Value = 'Morris'
$ws.Range("B7").Value = 'Ring Perpendicular to ND'
$ws.Range("B8").Value = 'Ring Perpendicular to RD'
$ws.Range("B9").Value = 'Ring Perpendicular to TD'
$ws.Range("B10").Value = 'Gaussian-Quadrature'
$ws.Range("B11").Value = 'Spiral-90deg-10rot-5space'
$ws.Range("B12").Value = 'Spiral-90deg-15rot-5space'
$ws.Range("B13").Value = 'Spiral-90deg-10rot-3space'
$ws.Range("B14").Value = 'NoRotation-tilt60deg'
$ws.Range("B15").Value = 'Rotation-NoTilt'
$ws.Range("B16").Value = 'Rotation-60detTilt'
$ws.Range("B17").Value = 'HexGrid-90degTilt5degRes'
$ws.Range("B18").Value = 'HexGrid-90degTilt22p5degRes'
$ws.Range("B19").Value = 'HexGrid-60degTilt5degRes'

# --- 2) Seed shared-string header row (C2:M2) left-to-right ---
$ws.Range("C2").Value = '[1, 1, 0]'
$ws.Range("D2").Value = '[2, 0, 0]'
$ws.Range("E2").Value = '[2, 1, 1]'
$ws.Range("F2").Value = '[2, 2, 0]'
$ws.Range("G2").Value = '[3, 1, 0]'
$ws.Range("H2").Value = '[2, 2, 2]'
$ws.Range("I2").Value = '[3, 2, 1]'
$ws.Range("J2").Value = '[4, 0, 0]'
$ws.Range("K2").Value = '2Pairs'
$ws.Range("L2").Value = '4Pairs'
$ws.Range("M2").Value = 'MaxUnique'

# --- 3) Make sure new rows 17:19 exist with the same formatting as the existing table rows ---
# (copy number/font/border formatting from row 16's A:M range down to rows 17:19)
$ws.Range("A16:M16").Copy() | Out-Null
$ws.Range("A17:M19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 4) Row 1 (numeric column index header) ---
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11

# --- 5) Column A (row index numbers), rows 2:19 ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17

# --- 6) Numeric data cells C:M for rows 3:19 ---
# Row 3
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 2.69
$ws.Range("E3").Value = 0.31
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 2.34
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.16
$ws.Range("J3").Value = 2.69
$ws.Range("K3").Value = 1.5
$ws.Range("L3").Value = 0.75
$ws.Range("M3").Value = 0.9166666666666666
# Row 4
$ws.Range("C4").Value = 1.3
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.68
$ws.Range("F4").Value = 1.3
$ws.Range("G4").Value = 0.01
$ws.Range("H4").Value = 0.87
$ws.Range("I4").Value = 1.7
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.34
$ws.Range("L4").Value = 0.8200000000000001
$ws.Range("M4").Value = 0.7599999999999999
# Row 5
$ws.Range("C5").Value = 7.59
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.03
$ws.Range("F5").Value = 7.59
$ws.Range("G5").Value = 0.06
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.63
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.015
$ws.Range("L5").Value = 3.8025
$ws.Range("M5").Value = 1.385
# Row 6
$ws.Range("C6").Value = 0.01
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 2.09
$ws.Range("F6").Value = 0.01
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 5.56
$ws.Range("I6").Value = 1.01
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1.045
$ws.Range("L6").Value = 0.5274999999999999
$ws.Range("M6").Value = 1.445
# Row 7
$ws.Range("C7").Value = 1.071780821917808
$ws.Range("D7").Value = 1.524520547945206
$ws.Range("E7").Value = 0.6320547945205479
$ws.Range("F7").Value = 1.071780821917808
$ws.Range("G7").Value = 1.505068493150685
$ws.Range("H7").Value = 0.119041095890411
$ws.Range("I7").Value = 0.7384931506849315
$ws.Range("J7").Value = 1.524520547945206
$ws.Range("K7").Value = 1.078287671232877
$ws.Range("L7").Value = 1.075034246575342
$ws.Range("M7").Value = 0.931826484018265
# Row 8
$ws.Range("C8").Value = 1.281052631578947
$ws.Range("D8").Value = 0.4068421052631579
$ws.Range("E8").Value = 0.8247368421052632
$ws.Range("F8").Value = 1.281052631578947
$ws.Range("G8").Value = 0.8621052631578947
$ws.Range("H8").Value = 1.455789473684211
$ws.Range("I8").Value = 1.041578947368421
$ws.Range("J8").Value = 0.4068421052631579
$ws.Range("K8").Value = 0.6157894736842106
$ws.Range("L8").Value = 0.9484210526315791
$ws.Range("M8").Value = 0.9786842105263158
# Row 9
$ws.Range("C9").Value = 0.8142105263157895
$ws.Range("D9").Value = 1.690526315789474
$ws.Range("E9").Value = 1.123157894736842
$ws.Range("F9").Value = 0.8142105263157895
$ws.Range("G9").Value = 0.6352631578947369
$ws.Range("H9").Value = 2.461578947368421
$ws.Range("I9").Value = 0.86
$ws.Range("J9").Value = 1.690526315789474
$ws.Range("K9").Value = 1.406842105263158
$ws.Range("L9").Value = 1.110526315789474
$ws.Range("M9").Value = 1.264122807017544
# Row 10
$ws.Range("C10").Value = 1.66636463905608
$ws.Range("D10").Value = 2.284038275547505
$ws.Range("E10").Value = 0.4213586778482891
$ws.Range("F10").Value = 1.66636463905608
$ws.Range("G10").Value = 1.319409511752389
$ws.Range("H10").Value = 0.7193734977715674
$ws.Range("I10").Value = 0.7084720708245341
$ws.Range("J10").Value = 2.284038275547505
$ws.Range("K10").Value = 1.352698476697897
$ws.Range("L10").Value = 1.509531557876988
$ws.Range("M10").Value = 1.186502778800061
# Row 11
$ws.Range("C11").Value = 0.8694892747283528
$ws.Range("D11").Value = 1.616730401604222
$ws.Range("E11").Value = 1.193267150561191
$ws.Range("F11").Value = 0.8694892747283528
$ws.Range("G11").Value = 0.5242297349077162
$ws.Range("H11").Value = 2.698049367619953
$ws.Range("I11").Value = 0.8927328530385661
$ws.Range("J11").Value = 1.616730401604222
$ws.Range("K11").Value = 1.404998776082707
$ws.Range("L11").Value = 1.13724402540553
$ws.Range("M11").Value = 1.29908313041
# Row 12
$ws.Range("C12").Value = 0.8675595687726367
$ws.Range("D12").Value = 1.626252271051958
$ws.Range("E12").Value = 1.195316210947415
$ws.Range("F12").Value = 0.8675595687726367
$ws.Range("G12").Value = 0.5247319076382049
$ws.Range("H12").Value = 2.701506108949616
$ws.Range("I12").Value = 0.8910329382874577
$ws.Range("J12").Value = 1.626252271051958
$ws.Range("K12").Value = 1.410784240999686
$ws.Range("L12").Value = 1.139171904886161
$ws.Range("M12").Value = 1.301066500941215
# Row 13
$ws.Range("C13").Value = 0.868348594528921
$ws.Range("D13").Value = 1.623752343166865
$ws.Range("E13").Value = 1.194378760251365
$ws.Range("F13").Value = 0.868348594528921
$ws.Range("G13").Value = 0.5234808631428263
$ws.Range("H13").Value = 2.697386664723736
$ws.Range("I13").Value = 0.8928245424531979
$ws.Range("J13").Value = 1.623752343166865
$ws.Range("K13").Value = 1.409065551709115
$ws.Range("L13").Value = 1.138707073119018
$ws.Range("M13").Value = 1.300028628044485
# Row 14
$ws.Range("C14").Value = 0.03339599999999845
$ws.Range("D14").Value = 2.626315999999997
$ws.Range("E14").Value = 0.8468600000000013
$ws.Range("F14").Value = 0.03339599999999845
$ws.Range("G14").Value = 1.617516
$ws.Range("H14").Value = 1.084823999999998
$ws.Range("I14").Value = 0.4687520000000004
$ws.Range("J14").Value = 2.626315999999997
$ws.Range("K14").Value = 1.736587999999999
$ws.Range("L14").Value = 0.8849919999999987
$ws.Range("M14").Value = 1.112943999999999
# Row 15
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 2.688625000000001
$ws.Range("E15").Value = 0.3101375
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 2.331137499999994
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0.16
$ws.Range("J15").Value = 2.688625000000001
$ws.Range("K15").Value = 1.499381250000001
$ws.Range("L15").Value = 0.7496906250000004
$ws.Range("M15").Value = 0.9149833333333325
# Row 16
$ws.Range("C16").Value = 0.4033276195839986
$ws.Range("D16").Value = 2.045043471462396
$ws.Range("E16").Value = 0.6026573177855999
$ws.Range("F16").Value = 0.4033276195839986
$ws.Range("G16").Value = 1.779536449331192
$ws.Range("H16").Value = 0.4065127591936034
$ws.Range("I16").Value = 0.5169440206848016
$ws.Range("J16").Value = 2.045043471462396
$ws.Range("K16").Value = 1.323850394623998
$ws.Range("L16").Value = 0.8635890071039984
$ws.Range("M16").Value = 0.9590036063402655
# Row 17
$ws.Range("C17").Value = 0.9894216124101333
$ws.Range("D17").Value = 0.9939417127373235
$ws.Range("E17").Value = 0.9907191400909077
$ws.Range("F17").Value = 0.9894216124101333
$ws.Range("G17").Value = 1.001245482278839
$ws.Range("H17").Value = 0.9869078065318353
$ws.Range("I17").Value = 0.9923555296373245
$ws.Range("J17").Value = 0.9939417127373235
$ws.Range("K17").Value = 0.9923304264141155
$ws.Range("L17").Value = 0.9908760194121244
$ws.Range("M17").Value = 0.992431880614394
# Row 18
$ws.Range("C18").Value = 1.190070044827548
$ws.Range("D18").Value = 0.7365972113002464
$ws.Range("E18").Value = 1.014202346751854
$ws.Range("F18").Value = 1.190070044827548
$ws.Range("G18").Value = 0.8909630636550357
$ws.Range("H18").Value = 1.11238213258114
$ws.Range("I18").Value = 0.9586304662879314
$ws.Range("J18").Value = 0.7365972113002464
$ws.Range("K18").Value = 0.8753997790260503
$ws.Range("L18").Value = 1.032734911926799
$ws.Range("M18").Value = 0.9838075442339593
# Row 19
$ws.Range("C19").Value = 1.074911693225758
$ws.Range("D19").Value = 0.7052513677721831
$ws.Range("E19").Value = 1.063920164590017
$ws.Range("F19").Value = 1.074911693225758
$ws.Range("G19").Value = 0.8232938332474752
$ws.Range("H19").Value = 1.184363173749573
$ws.Range("I19").Value = 1.069344838141505
$ws.Range("J19").Value = 0.7052513677721831
$ws.Range("K19").Value = 0.8845857661811001
$ws.Range("L19").Value = 0.9797487297034289
$ws.Range("M19").Value = 0.9868475117877518
